# Applies the FFXIV "Goblin Profits" leve-profit market-price refresh.
# For each affected row, currentAveragePrice(NQ/HQ) columns (H/I/J) were
# re-pulled from the market board, which ripples into the computed
# LevePriceNQ/HQ (K/L) and LeveProfitNQ/HQ (M/N) columns.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2392.0625
$ws.Range("I80").Value = 1242.2858
$ws.Range("K80").Value = 3726.8574
$ws.Range("M80").Value = -2728.8574
$ws.Range("H83").Value = 2392.0625
$ws.Range("I83").Value = 1242.2858
$ws.Range("K83").Value = 11180.5722
$ws.Range("M83").Value = -6188.572200000001
$ws.Range("H86").Value = 7203.294
$ws.Range("I86").Value = 6085.4
$ws.Range("K86").Value = 6085.4
$ws.Range("M86").Value = -4962.4
$ws.Range("H88").Value = 4851.353
$ws.Range("I88").Value = 2089.4285
$ws.Range("J88").Value = 6784.7
$ws.Range("K88").Value = 2089.4285
$ws.Range("L88").Value = 6784.7
$ws.Range("M88").Value = -1683.4285
$ws.Range("N88").Value = -7596.7
$ws.Range("H89").Value = 7203.294
$ws.Range("I89").Value = 6085.4
$ws.Range("K89").Value = 30427
$ws.Range("M89").Value = -24811
$ws.Range("H91").Value = 4851.353
$ws.Range("I91").Value = 2089.4285
$ws.Range("J91").Value = 6784.7
$ws.Range("K91").Value = 2089.4285
$ws.Range("L91").Value = 6784.7
$ws.Range("M91").Value = -685.4285
$ws.Range("N91").Value = -9592.700000000001
$ws.Range("H111").Value = 3545.6
$ws.Range("I111").Value = 3114.5
$ws.Range("J111").Value = 3833
$ws.Range("K111").Value = 9343.5
$ws.Range("L111").Value = 11499
$ws.Range("M111").Value = -6276.5
$ws.Range("N111").Value = -17633
$ws.Range("H116").Value = 9226.888999999999
$ws.Range("I116").Value = 8318.611000000001
$ws.Range("J116").Value = 11043.444
$ws.Range("K116").Value = 8318.611000000001
$ws.Range("L116").Value = 11043.444
$ws.Range("M116").Value = -4876.611000000001
$ws.Range("N116").Value = -17927.444
$ws.Range("H129").Value = 1986.1666
$ws.Range("I129").Value = 1600
$ws.Range("J129").Value = 2758.5
$ws.Range("K129").Value = 4800
$ws.Range("L129").Value = 8275.5
$ws.Range("M129").Value = 200
$ws.Range("N129").Value = -18275.5
$ws.Range("H131").Value = 5604.8887
$ws.Range("I131").Value = 2509.4
$ws.Range("J131").Value = 9474.25
$ws.Range("K131").Value = 7528.200000000001
$ws.Range("L131").Value = 28422.75
$ws.Range("M131").Value = -2488.200000000001
$ws.Range("N131").Value = -38502.75
$ws.Range("H135").Value = 835.4286
$ws.Range("I135").Value = 835.4286
$ws.Range("K135").Value = 7518.8574
$ws.Range("M135").Value = -4983.8574
$ws.Range("H137").Value = 2029.2084
$ws.Range("I137").Value = 1422.2142
$ws.Range("K137").Value = 4266.642599999999
$ws.Range("M137").Value = -1716.642599999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 139733.58
$ws.Range("I32").Value = 149710.83
$ws.Range("K32").Value = 149710.83
$ws.Range("M32").Value = -149423.83
$ws.Range("H45").Value = 2074.1667
$ws.Range("I45").Value = 1508.8235
$ws.Range("K45").Value = 1508.8235
$ws.Range("M45").Value = -1131.8235
$ws.Range("H102").Value = 7557.048
$ws.Range("I102").Value = 5724.9165
$ws.Range("K102").Value = 5724.9165
$ws.Range("M102").Value = -4102.9165
$ws.Range("H132").Value = 8312.647000000001
$ws.Range("I132").Value = 8802.532999999999
$ws.Range("K132").Value = 26407.599
$ws.Range("M132").Value = -23877.599

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2244.2666
$ws.Range("I94").Value = 2244.2666
$ws.Range("K94").Value = 2244.2666
$ws.Range("M94").Value = -1793.2666
$ws.Range("H107").Value = 4299.95
$ws.Range("I107").Value = 3167.7097
$ws.Range("K107").Value = 3167.7097
$ws.Range("M107").Value = -1247.7097
$ws.Range("H134").Value = 2265.5938
$ws.Range("I134").Value = 2241.3447
$ws.Range("K134").Value = 6724.034100000001
$ws.Range("M134").Value = -4189.034100000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3007.9583
$ws.Range("I58").Value = 3082.238
$ws.Range("J58").Value = 2488
$ws.Range("K58").Value = 3082.238
$ws.Range("L58").Value = 2488
$ws.Range("M58").Value = -2879.238
$ws.Range("N58").Value = -2894
$ws.Range("H94").Value = 1060.2727
$ws.Range("I94").Value = 815.25
$ws.Range("J94").Value = 1200.2858
$ws.Range("K94").Value = 815.25
$ws.Range("L94").Value = 1200.2858
$ws.Range("M94").Value = -364.25
$ws.Range("N94").Value = -2102.2858
$ws.Range("H99").Value = 3194.2856
$ws.Range("I99").Value = 2590.25
$ws.Range("K99").Value = 2590.25
$ws.Range("M99").Value = -1092.25
$ws.Range("H126").Value = 3194.2856
$ws.Range("I126").Value = 2590.25
$ws.Range("K126").Value = 7770.75
$ws.Range("M126").Value = -5300.75
$ws.Range("H134").Value = 31918.363
$ws.Range("I134").Value = 38489.277
$ws.Range("J134").Value = 2349.25
$ws.Range("K134").Value = 115467.831
$ws.Range("L134").Value = 7047.75
$ws.Range("M134").Value = -112932.831
$ws.Range("N134").Value = -12117.75
$ws.Range("H136").Value = 3007.9583
$ws.Range("I136").Value = 3082.238
$ws.Range("J136").Value = 2488
$ws.Range("K136").Value = 9246.714
$ws.Range("L136").Value = 7464
$ws.Range("M136").Value = -6696.714
$ws.Range("N136").Value = -12564

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 798.42426
$ws.Range("J34").Value = 557
$ws.Range("L34").Value = 1671
$ws.Range("N34").Value = -1839
$ws.Range("H55").Value = 2188523.5
$ws.Range("J55").Value = 3368500
$ws.Range("L55").Value = 10105500
$ws.Range("N55").Value = -10105854
$ws.Range("H140").Value = 2741.875
$ws.Range("I140").Value = 2675.7144
$ws.Range("K140").Value = 8027.1432
$ws.Range("M140").Value = -2847.1432

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 1000000000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H88").Value = 1000000000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H90").Value = 1000000000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H91").Value = 1000000000
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H107").Value = 2469.077
$ws.Range("I107").Value = 1156.8572
$ws.Range("K107").Value = 1156.8572
$ws.Range("M107").Value = 763.1428000000001
$ws.Range("H132").Value = 2770.68
$ws.Range("I132").Value = 2659.923
$ws.Range("K132").Value = 7979.768999999999
$ws.Range("M132").Value = -5449.768999999999
$ws.Range("H133").Value = 225934.19
$ws.Range("J133").Value = 225934.19
$ws.Range("L133").Value = 225934.19
$ws.Range("N133").Value = -236054.19

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2724.25
$ws.Range("I7").Value = 1799
$ws.Range("K7").Value = 1799
$ws.Range("M7").Value = -1687
$ws.Range("H22").Value = 3281.4443
$ws.Range("J22").Value = 4299.273
$ws.Range("L22").Value = 4299.273
$ws.Range("N22").Value = -4889.273
$ws.Range("H27").Value = 3281.4443
$ws.Range("J27").Value = 4299.273
$ws.Range("L27").Value = 4299.273
$ws.Range("N27").Value = -4513.273
$ws.Range("H126").Value = 2724.25
$ws.Range("I126").Value = 1799
$ws.Range("K126").Value = 5397
$ws.Range("M126").Value = -2927
$ws.Range("H132").Value = 3458.0344
$ws.Range("I132").Value = 2741.8235
$ws.Range("K132").Value = 8225.470499999999
$ws.Range("M132").Value = -5695.470499999999
$ws.Range("H136").Value = 19660.893
$ws.Range("I136").Value = 3247.6365
$ws.Range("K136").Value = 9742.9095
$ws.Range("M136").Value = -7192.9095

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 923.3333
$ws.Range("I100").Value = 515
$ws.Range("K100").Value = 1030
$ws.Range("M100").Value = -489
$ws.Range("H126").Value = 2881.6843
$ws.Range("I126").Value = 2278.5625
$ws.Range("J126").Value = 6098.3335
$ws.Range("K126").Value = 6835.6875
$ws.Range("L126").Value = 18295.0005
$ws.Range("M126").Value = -4365.6875
$ws.Range("N126").Value = -23235.0005
$ws.Range("H132").Value = 6011.381
$ws.Range("I132").Value = 6092.26
$ws.Range("K132").Value = 18276.78
$ws.Range("M132").Value = -15746.78
$ws.Range("H136").Value = 2825.3447
$ws.Range("I136").Value = 957.1
$ws.Range("K136").Value = 2871.3
$ws.Range("M136").Value = -321.3000000000002
